$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as text so values like "1.000" or "29.211.59"
# are not auto-converted to numbers by Excel, matching the original inlineStr cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.211.59'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.853.09'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '0.6978'
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").Value = '237.70'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '0.07875'
$ws.Range("E8").Value = '  +1.34%  '
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("D10").Value = '23.78'
$ws.Range("E10").Value = '  +2.18%  '
$ws.Range("D11").Value = '0.08092'
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").Value = '1.854.65'
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").Value = '5.182'
$ws.Range("E13").Value = '  -0.07%  '
$ws.Range("D14").Value = '0.7045'
$ws.Range("D15").Value = '89.48'
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").Value = '29.249.18'
$ws.Range("D17").Value = '5.805'
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").Value = '0.000007810'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Value = '13.20'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("D20").Value = '235.48'
$ws.Range("E20").Value = '  +0.25%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '2.107.50'
$ws.Range("E22").Value = '  -0.60%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '7.506'
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").Value = '162.27'
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").Value = '8.853'
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("D27").Value = '0.1418'
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").Value = '18.02'
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").Value = '1.918'
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("D30").Value = '1.405'
$ws.Range("E30").Value = '  +0.38%  '
$ws.Range("D31").Value = '1.478'
$ws.Range("E31").Value = '  -0.45%  '
$ws.Range("D32").Value = '4.322'
$ws.Range("E32").Value = '  -4.73%  '
$ws.Range("D33").Value = '4.011'
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '0.05158'
$ws.Range("E34").Value = '  -0.89%  '
$ws.Range("E35").Value = '  -1.98%  '
$ws.Range("D36").Value = '0.7107'
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").Value = '0.9976'
$ws.Range("E37").Value = '  -2.57%  '
$ws.Range("D38").Value = '2.680'
$ws.Range("D39").Value = '0.01843'
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").Value = '2.709'
$ws.Range("E40").Value = '  +1.05%  '
$ws.Range("D41").Value = '1.152.53'
$ws.Range("E41").Value = '  +5.50%  '
$ws.Range("D42").Value = '0.9227'
$ws.Range("E42").Value = '  -0.75%  '
$ws.Range("D43").Value = '5.972'
$ws.Range("E43").Value = '  -0.66%  '
$ws.Range("D44").Value = '0.4228'
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("D45").Value = '69.97'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '102.95'
$ws.Range("D48").Value = '0.5300'
$ws.Range("E48").Value = '  -2.73%  '
$ws.Range("D49").Value = '1.736'
$ws.Range("E49").Value = '  -3.29%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.113'
$ws.Range("E50").Value = '  -0.50%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '6.945'
$ws.Range("E51").Value = '  -0.88%  '
